$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.951.79"
$ws.Range("E2").Value = "  -1.29%  "

$ws.Range("D3").Value = "1.993.46"
$ws.Range("E3").Value = "  -3.11%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.81"
$ws.Range("E5").Value = "  -2.91%  "

$ws.Range("E6").Value = "  -2.63%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.51"
$ws.Range("E8").Value = "  -4.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.375"
$ws.Range("E9").Value = "  -2.53%  "

$ws.Range("E10").Value = "  +0.91%  "

$ws.Range("E11").Value = "  -3.73%  "

$ws.Range("D12").Value = "2.285.51"
$ws.Range("E12").Value = "  -3.16%  "

$ws.Range("E13").Value = "  -4.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.05"
$ws.Range("E14").Value = "  -4.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.736"
$ws.Range("E15").Value = "  -3.17%  "

$ws.Range("E16").Value = "  -4.50%  "

$ws.Range("D17").Value = "1.999.23"
$ws.Range("E17").Value = "  -2.79%  "

$ws.Range("D18").Value = "36.786.53"
$ws.Range("E18").Value = "  -1.74%  "

$ws.Range("E19").Value = "  -0.58%  "

$ws.Range("E20").Value = "  -1.77%  "

$ws.Range("D21").Value = "0.0₃0812"
$ws.Range("E21").Value = "  -1.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "222.14"
$ws.Range("E22").Value = "  -1.99%  "

$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("E24").Value = "  +0.60%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.17"
$ws.Range("E25").Value = "  -7.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.58"
$ws.Range("E26").Value = "  -2.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.08"
$ws.Range("E27").Value = "  -8.05%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.59"
$ws.Range("E28").Value = "  -3.25%  "

$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.34"
$ws.Range("E29").Value = "  -0.48%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.123"
$ws.Range("E30").Value = "  -5.14%  "

$ws.Range("E31").Value = "  -4.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.48"
$ws.Range("E32").Value = "  -1.27%  "

$ws.Range("E33").Value = "  -2.64%  "

$ws.Range("E34").Value = "  -4.51%  "

$ws.Range("E35").Value = "  -7.16%  "

$ws.Range("E36").Value = "  +1.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.13"
$ws.Range("E38").Value = "  -4.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.26"
$ws.Range("E39").Value = "  -1.03%  "

$ws.Range("D40").Value = "1.464.66"
$ws.Range("E40").Value = "  -0.84%  "

$ws.Range("E41").Value = "  -4.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "94.36"
$ws.Range("E42").Value = "  -4.17%  "

$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0910"
$ws.Range("E43").Value = "  -4.97%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.26"
$ws.Range("E44").Value = "  -2.17%  "

$ws.Range("E45").Value = "  -4.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.12"
$ws.Range("E46").Value = "  -6.69%  "

$ws.Range("E47").Value = "  -2.93%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.13"
$ws.Range("E48").Value = "  -1.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.89"
$ws.Range("E49").Value = "  -2.21%  "

$ws.Range("D50").Value = "2.176.70"
$ws.Range("E50").Value = "  -3.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.12"
$ws.Range("E51").Value = "  -3.65%  "
